$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated loading_percent results for the case with 380 kV (rows 2-25, cols C:H and J:O)
$loadingPercent = @{
    2 = @(13.46550996381809;9.037691232435078;14.17234311566548;35.06066457061306;37.84169845501361;16.63956290156992;10.60591078380693;14.57093058763786;9.734168444495468;18.0610291667824;18.61425711262354;26.52479751314108)
    3 = @(13.44673030312472;9.02616321257651;14.19491578177187;35.13862235457602;37.92907678307952;16.68980471128704;10.62926462510353;14.22938384432856;9.749507728923547;17.9388952954114;18.65181194000182;26.60548198698259)
    4 = @(13.43780655411631;9.020235064787107;14.210580812911;35.19401036964008;37.99306399821081;16.72321316538785;10.64449783093952;14.01679174851476;9.759507345446798;17.86560445839306;18.67660487773862;26.66026848123292)
    5 = @(13.4348290307675;9.018110558995367;14.21741882707776;35.21846970849921;38.02173049634263;16.73747101671368;10.65093073946713;13.92955098204514;9.763728810624576;17.836188320585;18.68714504118716;26.68391153685605)
    6 = @(13.43437450648271;9.017775435274789;14.21858172966471;35.22264511239875;38.02664679753678;16.73987739475895;10.652012538143;13.91503152525671;9.764438643689745;17.83133166047465;18.68892163448967;26.6879169381837)
    7 = @(13.43776372585182;9.020205231029887;14.21067119242767;35.19433259566969;37.99344012423033;16.72340284576795;10.64458367470194;14.01561749177645;9.759563683754033;17.86520588820844;18.67674525628028;26.66058200953805)
    8 = @(13.45849583214949;9.033478944338565;14.17975174020582;35.08598154344475;37.86967656093017;16.65635497248464;10.61377796213094;14.45383155415379;9.739337025918239;18.01857861514373;18.62684656714368;26.55152755182602)
    9 = @(13.51966770612579;9.068538540689541;14.13342389840495;34.93331159613591;37.70931720377756;16.54518675884836;10.56043835691506;15.28546803221992;9.704267094812041;18.33179072211053;18.54272052185889;26.37938288026921)
    10 = @(13.57686290844856;9.099665908205573;14.10808053890852;34.85775603028713;37.6420900833538;16.47589665039495;10.52552946466915;15.87354015062232;9.681278124975936;18.56807224324731;18.48923407545044;26.27845052172567)
    11 = @(13.60547899936383;9.114960517480947;14.09843242597135;34.83135765995134;37.62255614245513;16.44706455384445;10.51057136927699;16.13494362035378;9.671417746524547;18.6765903000272;18.46669894372146;26.23810248709434)
    12 = @(13.61668255211009;9.120912386634542;14.09504875491338;34.82250885796949;37.616751415526;16.43653323690906;10.50503924647758;16.23296486150239;9.667769398265589;18.71780513275917;18.45842302743779;26.22362587813762)
    13 = @(13.61425343736812;9.119623472814705;14.09576549672553;34.82436354664802;37.6179307001382;16.43878413962586;10.50622481582665;16.21189854166176;9.668551335362492;18.70892381350299;18.46019394499641;26.22670796371054)
    14 = @(13.60639339952755;9.115446989858579;14.09814864443319;34.83060665662953;37.6220466537988;16.44619038348361;10.51011359143339;16.14302773503331;9.671115881871524;18.67997879137286;18.46601291901005;26.23689539805815)
    15 = @(13.60162653338592;9.112909535314332;14.0996435168148;34.83458023470534;37.62477525595872;16.4507772965061;10.51251278061264;16.10071402207499;9.67269787330684;18.66226413481901;18.46961074359691;26.24324003348554)
    16 = @(13.57504452483751;9.098688945613086;14.10874886162635;34.85964176884217;37.64358932398446;16.47783500846536;10.52652553119003;15.8563261653142;9.681934515927566;18.56099879695131;18.49074288764234;26.2811995357236)
    17 = @(13.55939832720034;9.090253569988556;14.11481598995891;34.87705915842604;37.65796374965495;16.49512272328776;10.53535778114395;15.70477189024367;9.687753663850499;18.49912100042807;18.50416635048109;26.30591355603024)
    18 = @(13.55064396154423;9.085508750267483;14.11848270076518;34.88782751750631;37.6672712781201;16.50531918594859;10.54052467809243;15.61703123781927;9.691156936184186;18.4636293325733;18.5120562676233;26.32065225951002)
    19 = @(13.5477221321173;9.083920700455845;14.11975461311491;34.89160232290532;37.67060109452482;16.50881498445981;10.54228902494503;15.58722853069636;9.692318897684189;18.45163025607412;18.51475671989093;26.32573244508559)
    20 = @(13.56103858817022;9.091140479014397;14.11415181292411;34.87512737960311;37.65632592848473;16.49325622780885;10.53440859052309;15.72096481634041;9.687128386222158;18.50569796733795;18.50271990371213;26.30322847579107)
    21 = @(13.60869216812627;9.116669402981227;14.09744133714907;34.82874174807261;37.62079445858077;16.44400449173246;10.50896777966551;16.16328365595188;9.670360292973962;18.68847756731015;18.46429675773764;26.23388131296989)
    22 = @(13.64197398680337;9.134285905351796;14.08809272092637;34.80511591666671;37.60685497406112;16.41407012416411;10.4931110139757;16.44669489099936;9.659899943808746;18.80862962476244;18.44068650191923;26.19323620589794)
    23 = @(13.62401746712107;9.124799414685478;14.09293856126833;34.81711304341247;37.61344446913027;16.42984031145486;10.50150372497482;16.29597887824067;9.665437323117372;18.74444757163776;18.45315055867073;26.21450071886217)
    24 = @(13.56029627553097;9.090739180817124;14.11445153088482;34.87599838553926;37.65706313734943;16.49409926824515;10.53483744197297;15.71364588798496;9.687410894274842;18.50272426092076;18.503373304638;26.30444074889567)
    25 = @(13.50094713702938;9.058100700967525;14.1444279756941;34.96819288483857;37.74383937167968;16.57308587606429;10.57411438443765;15.06410334760292;9.713265063196705;18.24587114410379;18.56401430478504;26.42147491408848)
}

$cols = @(3,4,5,6,7,8,10,11,12,13,14,15)  # C,D,E,F,G,H,J,K,L,M,N,O (column I stays 0)

foreach ($row in $loadingPercent.Keys) {
    $values = $loadingPercent[$row]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $ws.Cells.Item($row, $cols[$i]).Value = $values[$i]
    }
}
